# Weekly crime-data refresh: shift report week forward, update volume/
# number, and refresh the Week-to-Date / 28-Day / Year-to-Date / 2-Year
# statistics for the affected precinct rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume number + report week dates -----------------
# A8 currently reads "Volume 31   Number  34" -> bump issue number to 35
$volCell = $ws.Cells.Item(8, 1)
$volText = $volCell.Characters(1, 300).Text
$numPos = $volText.Length - 1
$volCell.Characters($numPos, 2).Text = "35"

# C9 currently reads "Report Covering the Week  8/19/2024  Through  8/25/2024"
# -> shift the week forward by one week
$weekCell = $ws.Cells.Item(9, 3)
$weekCell.Characters(27, 9).Text = "8/26/2024"
$weekCell.Characters(47, 9).Text = "9/1/2024"

# --- Precinct statistics table (rows 14-33) --------------------------
$ws.Cells.Item(14, 14).Value = -80

$ws.Cells.Item(15, 14).Value = -52

$ws.Cells.Item(16, 3).Value = 6
$ws.Cells.Item(16, 4).Value = "0"
$ws.Cells.Item(16, 5).Value = "***.*"
$ws.Cells.Item(16, 7).Value = 13
$ws.Cells.Item(16, 8).Value = 38.461538461538
$ws.Cells.Item(16, 9).Value = 138
$ws.Cells.Item(16, 11).Value = 4.545454545454
$ws.Cells.Item(16, 12).Value = -6.756756756756
$ws.Cells.Item(16, 13).Value = -15.337423312883
$ws.Cells.Item(16, 14).Value = -75.618374558303

$ws.Cells.Item(17, 3).Value = 6
$ws.Cells.Item(17, 4).Value = "0"
$ws.Cells.Item(17, 5).Value = "***.*"
$ws.Cells.Item(17, 6).Value = 27
$ws.Cells.Item(17, 7).Value = 19
$ws.Cells.Item(17, 8).Value = 42.105263157894
$ws.Cells.Item(17, 9).Value = 266
$ws.Cells.Item(17, 11).Value = 12.711864406779
$ws.Cells.Item(17, 12).Value = 11.764705882352
$ws.Cells.Item(17, 13).Value = 189.130434782609
$ws.Cells.Item(17, 14).Value = -2.919708029197

$ws.Cells.Item(18, 3).Value = "0"
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 5).Value = -100
$ws.Cells.Item(18, 6).Value = 10
$ws.Cells.Item(18, 7).Value = 14
$ws.Cells.Item(18, 8).Value = -28.571428571428
$ws.Cells.Item(18, 10).Value = 98
$ws.Cells.Item(18, 11).Value = -14.285714285714
$ws.Cells.Item(18, 12).Value = -13.40206185567
$ws.Cells.Item(18, 13).Value = -58.20895522388
$ws.Cells.Item(18, 14).Value = -89.244558258642

$ws.Cells.Item(19, 3).Value = 16
$ws.Cells.Item(19, 4).Value = 10
$ws.Cells.Item(19, 5).Value = 60
$ws.Cells.Item(19, 6).Value = 51
$ws.Cells.Item(19, 7).Value = 49
$ws.Cells.Item(19, 8).Value = 4.081632653061
$ws.Cells.Item(19, 9).Value = 362
$ws.Cells.Item(19, 10).Value = 422
$ws.Cells.Item(19, 11).Value = -14.218009478673
$ws.Cells.Item(19, 12).Value = -20.264317180616
$ws.Cells.Item(19, 13).Value = 66.055045871559
$ws.Cells.Item(19, 14).Value = -5.483028720626

$ws.Cells.Item(20, 3).Value = 3
$ws.Cells.Item(20, 4).Value = 4
$ws.Cells.Item(20, 5).Value = -25
$ws.Cells.Item(20, 6).Value = 26
$ws.Cells.Item(20, 7).Value = 23
$ws.Cells.Item(20, 8).Value = 13.043478260869
$ws.Cells.Item(20, 9).Value = 217
$ws.Cells.Item(20, 10).Value = 167
$ws.Cells.Item(20, 11).Value = 29.940119760479
$ws.Cells.Item(20, 12).Value = 45.637583892617
$ws.Cells.Item(20, 13).Value = 14.210526315789
$ws.Cells.Item(20, 14).Value = -89.850327408793

$ws.Cells.Item(21, 3).Value = 31
$ws.Cells.Item(21, 4).Value = 16
$ws.Cells.Item(21, 5).Value = 93.75
$ws.Cells.Item(21, 6).Value = 132
$ws.Cells.Item(21, 7).Value = 119
$ws.Cells.Item(21, 8).Value = 10.924369747899
$ws.Cells.Item(21, 9).Value = 1082
$ws.Cells.Item(21, 10).Value = 1069
$ws.Cells.Item(21, 11).Value = 1.216089803554
$ws.Cells.Item(21, 12).Value = -2.0814479638
$ws.Cells.Item(21, 13).Value = 22.536806342015
$ws.Cells.Item(21, 14).Value = -74.127211860353

$ws.Cells.Item(22, 6).Value = 2
$ws.Cells.Item(22, 7).Value = 3
$ws.Cells.Item(22, 8).Value = -33.333333333333
$ws.Cells.Item(22, 10).Value = 16
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = -5.882352941176

$ws.Cells.Item(24, 3).Value = 38
$ws.Cells.Item(24, 4).Value = 16
$ws.Cells.Item(24, 5).Value = 137.5
$ws.Cells.Item(24, 6).Value = 130
$ws.Cells.Item(24, 7).Value = 82
$ws.Cells.Item(24, 8).Value = 58.536585365853
$ws.Cells.Item(24, 9).Value = 851
$ws.Cells.Item(24, 10).Value = 860
$ws.Cells.Item(24, 11).Value = -1.046511627906
$ws.Cells.Item(24, 12).Value = -10.042283298097
$ws.Cells.Item(24, 13).Value = 98.831775700934

$ws.Cells.Item(25, 3).Value = 22
$ws.Cells.Item(25, 4).Value = 7
$ws.Cells.Item(25, 5).Value = 214.285714285714
$ws.Cells.Item(25, 6).Value = 76
$ws.Cells.Item(25, 7).Value = 36
$ws.Cells.Item(25, 8).Value = 111.111111111111
$ws.Cells.Item(25, 9).Value = 445
$ws.Cells.Item(25, 10).Value = 380
$ws.Cells.Item(25, 11).Value = 17.105263157894
$ws.Cells.Item(25, 12).Value = 8.80195599022

$ws.Cells.Item(26, 3).Value = 17
$ws.Cells.Item(26, 4).Value = 7
$ws.Cells.Item(26, 5).Value = 142.857142857143
$ws.Cells.Item(26, 6).Value = 56
$ws.Cells.Item(26, 8).Value = 24.444444444444
$ws.Cells.Item(26, 9).Value = 432
$ws.Cells.Item(26, 10).Value = 371
$ws.Cells.Item(26, 11).Value = 16.44204851752
$ws.Cells.Item(26, 12).Value = 19.6675900277
$ws.Cells.Item(26, 13).Value = 23.428571428571

$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 6).Value = 3
$ws.Cells.Item(27, 7).Value = 5
$ws.Cells.Item(27, 8).Value = -40
$ws.Cells.Item(27, 9).Value = 20
$ws.Cells.Item(27, 10).Value = 22
$ws.Cells.Item(27, 11).Value = -9.090909090909
$ws.Cells.Item(27, 12).Value = -20

$ws.Cells.Item(28, 3).Value = 2
$ws.Cells.Item(28, 4).Value = 3
$ws.Cells.Item(28, 5).Value = -33.333333333333
$ws.Cells.Item(28, 6).Value = 6
$ws.Cells.Item(28, 7).Value = 5
$ws.Cells.Item(28, 8).Value = 20
$ws.Cells.Item(28, 9).Value = 31
$ws.Cells.Item(28, 10).Value = 36
$ws.Cells.Item(28, 11).Value = -13.888888888888
$ws.Cells.Item(28, 12).Value = -16.216216216216

$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 6).Value = 1
$ws.Cells.Item(33, 9).Value = 3
$ws.Cells.Item(33, 11).Value = -40
$ws.Cells.Item(33, 12).Value = 200

